$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

$data = @(
    @(0, 5351.066666666667, 5697, 4691, 0.1295172850290934),
    @(1, 5241.366666666667, 5584, 4762, 0.1293060064315796),
    @(2, 5180.066666666667, 5603, 4199, 0.1303920904795329),
    @(3, 5670.966666666666, 5963, 5325, 0.1302790562311808),
    @(4, 5042.566666666667, 5358, 4432, 0.1355497280756633),
    @(5, 5010.633333333333, 5413, 4378, 0.1309496641159058),
    @(6, 5529.1, 5984, 4872, 0.1349993944168091),
    @(7, 5383.833333333333, 5789, 4916, 0.1341822942097982),
    @(8, 5328.633333333333, 5577, 4668, 0.1346200386683146),
    @(9, 5333.666666666667, 5927, 4877, 0.129835311571757)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
